# Update the "Förändrad" (Changed) date column (C) for all data rows
# (rows 2 through 119) from serial date 45186 (2023-09-17) to
# serial date 45188 (2023-09-19).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 119
for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45186) {
        $cell.Value2 = 45188
    }
}
